$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for the "Contacts" (Kontakter) detail-page fields, appended
# after the existing data (rows 335-352).
$data = @(
    @(335, 'Edit', 24695, 'Kontakter', 'Namn'),
    @(336, 'Edit', 24694, 'Kontakter', 'Signatur'),
    @(337, 'Edit', 24696, 'Kontakter', 'Titel'),
    @(338, 'Edit', 24697, 'Kontakter', 'Telefon'),
    @(339, 'Edit', 24698, 'Kontakter', 'Telefon 2'),
    @(340, 'Edit', 24699, 'Kontakter', 'Mobiltelefon'),
    @(341, 'Edit', 24700, 'Kontakter', 'E-post'),
    @(342, 'Edit', 24701, 'Kontakter', 'Kundnummer'),
    @(343, 'Edit', 24702, 'Kontakter', 'Leverantörsnummer'),
    @(344, 'Edit', 24703, 'Kontakter', 'Företagsnamn'),
    @(345, 'Edit', 24704, 'Kontakter', 'Postadress'),
    @(346, 'Edit', 24705, 'Kontakter', 'Postadress 2'),
    @(347, 'Edit', 24707, 'Kontakter', 'Besöksadress'),
    @(348, 'Edit', 24708, 'Kontakter', 'Postnummer'),
    @(349, 'Edit', 24709, 'Kontakter', 'Ort'),
    @(350, 'Edit', 24711, 'Kontakter', 'Landskod'),
    @(351, 'Edit', 24710, 'Kontakter', 'Land'),
    @(352, 'Edit', 24706, 'Kontakter', 'Anteckningar')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Match the final view state recorded in the target workbook: scrolled
# down to the new rows, with the last new row's B:D cells selected.
$ws.Activate()
$ws.Range("B352:D352").Select()
